$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.629.21"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "'1.601.74"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'212.14"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'27.90"
$ws.Range("E8").Value = "  +6.16%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'1.831.20"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "'1.600.67"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "'0.544"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "'29.638.14"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "'63.86"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "'242.44"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'9.38"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "'155.06"
$ws.Range("D26").Value = "'15.43"
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "'6.42"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "'0.0482"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "'1.06"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("D34").Value = "'1.425.17"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("E36").Value = "  +5.48%  "
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("D38").Value = "'2.30"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("D40").Value = "'57.35"
$ws.Range("E40").Value = "  +8.29%  "
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("D42").Value = "'1.97"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "'0.0495"
$ws.Range("E43").Value = "  +4.55%  "
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.18"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'0.977"
$ws.Range("E47").Value = "  +16.88%  "
$ws.Range("D48").Value = "'5.34"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'1.741.17"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "'86.84"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").Value = "'0.0₆0102"
$ws.Range("E51").Value = "  +0.57%  "
